$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column A for "Lab. #" — shifts existing data (old columns A-H) to B-I.
$ws.Columns.Item(1).Insert()

# Header for the new column
$ws.Range("A1").Value = "Lab. #"

# Lab numbers for data rows 2-16 (the 9186 rows are the repeated QC/blank sample).
$labNums = @{
    2  = 9186
    3  = 9715
    4  = 9186
    5  = 9716
    6  = 9186
    7  = 9717
    8  = 9186
    9  = 9718
    10 = 9186
    11 = 9719
    12 = 9186
    13 = 9720
    14 = 9186
    15 = 9721
    16 = 9186
}

foreach ($r in 2..16) {
    $ws.Cells.Item($r, 1).Value = $labNums[$r]
}

# Highlight the repeated "9186" rows with a light-green solid fill (RGB D8E4BC).
$highlightRows = @(2, 4, 6, 8, 10, 12, 14, 16)
foreach ($r in $highlightRows) {
    $rowRange = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 9))
    $rowRange.Interior.Color = 12379352
}

# Restore/adjust column widths to match the new 9-column layout.
$ws.Columns.Item(1).ColumnWidth = 6.833333333333333
$ws.Columns.Item(2).ColumnWidth = 22.833333333333332
$ws.Columns.Item(3).ColumnWidth = 18.833333333333332
$ws.Columns.Item(4).ColumnWidth = 22.833333333333332
$ws.Columns.Item(5).ColumnWidth = 22.833333333333332
$ws.Columns.Item(6).ColumnWidth = 21.833333333333332
$ws.Columns.Item(7).ColumnWidth = 19.833333333333332
$ws.Columns.Item(8).ColumnWidth = 22.833333333333332
$ws.Columns.Item(9).ColumnWidth = 19.833333333333332
